$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '25.864.75'
$ws.Range("E2").Value = '  +0.07%  '
$ws.Range("D3").Value = '1.736.22'
$ws.Range("E3").Value = '  -0.77%  '
$ws.Range("D4").Value = '''0.9991'
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").Value = '''233.68'
$ws.Range("E5").Value = '  -1.08%  '
$ws.Range("D6").Value = '''0.9993'
$ws.Range("E6").Value = '  +0.01%  '
$ws.Range("D7").Value = '''0.5176'
$ws.Range("E7").Value = '  +0.88%  '
$ws.Range("D8").Value = '''0.2772'
$ws.Range("E8").Value = '  +3.96%  '
$ws.Range("D9").Value = '''39.32'
$ws.Range("E9").Value = '  -2.66%  '
$ws.Range("D10").Value = '''0.06117'
$ws.Range("E10").Value = '  -0.94%  '
$ws.Range("D11").Value = '1.739.10'
$ws.Range("E11").Value = '  -0.67%  '
$ws.Range("D12").Value = '''0.07050'
$ws.Range("E12").Value = '  +1.57%  '
$ws.Range("D13").Value = '''15.27'
$ws.Range("E13").Value = '  -0.74%  '
$ws.Range("D14").Value = '''0.6401'
$ws.Range("E14").Value = '  +2.06%  '
$ws.Range("D15").Value = '''4.508'
$ws.Range("E15").Value = '  +0.69%  '
$ws.Range("D16").Value = '''76.88'
$ws.Range("E16").Value = '  -1.36%  '
$ws.Range("D17").Value = '''0.9993'
$ws.Range("E17").Value = '  -0.03%  '
$ws.Range("D18").Value = '''0.9993'
$ws.Range("E18").Value = '  -0.04%  '
$ws.Range("D19").Value = '25.829.45'
$ws.Range("E19").Value = '  -0.16%  '
$ws.Range("E20").Value = '  -1.04%  '
$ws.Range("D21").Value = '''0.000006623'
$ws.Range("E21").Value = '  -0.40%  '
$ws.Range("D22").Value = '1.959.22'
$ws.Range("E22").Value = '  -2.27%  '
$ws.Range("D23").Value = '''4.134'
$ws.Range("E23").Value = '  +1.92%  '
$ws.Range("D24").Value = '''8.754'
$ws.Range("E24").Value = '  +6.01%  '
$ws.Range("D25").Value = '''5.136'
$ws.Range("E25").Value = '  -0.57%  '
$ws.Range("D26").Value = '''139.97'
$ws.Range("E26").Value = '  +2.70%  '
$ws.Range("D27").Value = '''1.512'
$ws.Range("E27").Value = '  +2.14%  '
$ws.Range("D28").Value = '''15.00'
$ws.Range("E28").Value = '  -0.87%  '
$ws.Range("D29").Value = '''1.786'
$ws.Range("E29").Value = '  +0.52%  '
$ws.Range("D30").Value = '''102.00'
$ws.Range("E30").Value = '  -0.77%  '
$ws.Range("D31").Value = '''0.08299'
$ws.Range("E31").Value = '  +0.33%  '
$ws.Range("D32").Value = '''3.683'
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("D33").Value = '''3.443'
$ws.Range("E33").Value = '  +1.27%  '
$ws.Range("D34").Value = '''0.04507'
$ws.Range("E34").Value = '  +2.81%  '
$ws.Range("D35").Value = '''2.617'
$ws.Range("E35").Value = '  -0.76%  '
$ws.Range("D36").Value = '''0.9780'
$ws.Range("E36").Value = '  -2.02%  '
$ws.Range("E37").Value = '  +1.19%  '
$ws.Range("D38").Value = '''2.664'
$ws.Range("E38").Value = '  -0.20%  '
$ws.Range("D39").Value = '''0.01583'
$ws.Range("E39").Value = '  +1.64%  '
$ws.Range("D40").Value = '''1.941'
$ws.Range("E40").Value = '  +0.45%  '
$ws.Range("D41").Value = '''0.9991'
$ws.Range("E41").Value = '  +0.02%  '
$ws.Range("D42").Value = '''100.51'
$ws.Range("E42").Value = '  -1.77%  '
$ws.Range("D43").Value = '''0.3834'
$ws.Range("D44").Value = '''0.7243'
$ws.Range("E44").Value = '  -3.12%  '
$ws.Range("D45").Value = '''4.977'
$ws.Range("D46").Value = '''0.05378'
$ws.Range("E46").Value = '  -2.02%  '
$ws.Range("D48").Value = '''6.251'
$ws.Range("E48").Value = '  +4.96%  '
$ws.Range("D49").Value = '''53.00'
$ws.Range("E49").Value = '  +0.89%  '
$ws.Range("E50").Value = '  -0.26%  '
$ws.Range("D51").Value = '''7.614'
$ws.Range("E51").Value = '  +1.91%  '
